# Generate Report for Handoff
# Regenerates the localization-status report: the two tracked source files
# have moved on to a new handoff cycle (new GUIDs / new content hash), the
# previous "handed back" state is cleared, and the now-irrelevant
# "Latest Target File" / "Latest Handback File" columns are cleared out
# because nothing has been handed back yet in this cycle.

$wb = $excel.ActiveWorkbook

# ---- New identifiers for this handoff run ----
$guid1 = "cc9c4d88-a1e4-4904-acb4-021401f05c23"
$guid2 = "ffff76ca6570-403e-4963-9d34-7be281117591"
$hash  = "6c0a929d12438973a48157a11a42e7268d01a887"

$md1 = "$guid1.md"
$md2 = "$guid2.md"
$xlfZh = "$guid1.$hash.zh-cn.xlf"
$xlfDe = "$guid1.$hash.de-de.xlf"

$status        = "Ready for handoff"
$handoffDtZh   = "2016-03-17 14:50:52"
$handoffDtDe   = "2016-03-17 14:50:56"
$handbackNone  = "0001-01-01 00:00:00"
$overviewDate  = "2016-50-17 14:50:56"

# ---- URLs (same repos/commit pattern as before, new filenames) ----
$srcRepoCommit = "d4731bbef80a452a3fa136653c559bc7ba8d16a5"
$urlMd1 = "https://github.com/OpenLocalizationTest/oltest/blob/$srcRepoCommit/e2e/$md1"
$urlMd2 = "https://github.com/OpenLocalizationTest/oltest/blob/$srcRepoCommit/e2e/$md2"

$handoffCommitZh = "a8abbb53f2c64c021fa370a9a834f5a55f3b6eba"
$handoffCommitDe = "bcd14aa45121d51d5299f46177d81b756362cbe2"
$urlXlfZh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$handoffCommitZh/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZh"
$urlXlfDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$handoffCommitDe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDe"

# =====================================================================
# Sheet "Overview"
# =====================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# Column B = zh-cn status, Column C = de-de status, Column D = Latest Handoff Date
$wsOverview.Range("B2").Value = $status
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C2").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D2").Value = $overviewDate
$wsOverview.Range("D3").Value = $overviewDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $urlMd1, "", "", $md1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $urlMd2, "", "", $md2)

# =====================================================================
# Sheet "zh-cn"
# =====================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

# "Latest Target File" (F) / "Latest Handback File" (G) no longer apply -
# nothing has been handed back yet in this cycle, so clear those cells.
$wsZh.Range("F2:G3").Clear()

$wsZh.Range("C2").Value = $status
$wsZh.Range("E2").Value = $handoffDtZh
$wsZh.Range("H2").Value = $handbackNone

$wsZh.Range("C3").Value = $status
$wsZh.Range("E3").Value = $handoffDtZh
$wsZh.Range("H3").Value = $handbackNone

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlMd1, "", "", $md1)
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $urlMd1, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $urlXlfZh, "", "", $xlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlMd2, "", "", $md2)
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $urlMd2, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $urlXlfZh, "", "", $xlfZh)

# =====================================================================
# Sheet "de-de"
# =====================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("F2:G3").Clear()

$wsDe.Range("C2").Value = $status
$wsDe.Range("E2").Value = $handoffDtDe
$wsDe.Range("H2").Value = $handbackNone

$wsDe.Range("C3").Value = $status
$wsDe.Range("E3").Value = $handoffDtDe
$wsDe.Range("H3").Value = $handbackNone

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlMd1, "", "", $md1)
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $urlMd1, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $urlXlfDe, "", "", $xlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlMd2, "", "", $md2)
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $urlMd2, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $urlXlfDe, "", "", $xlfDe)
